$d = $word.ActiveDocument

# --- Change 1: append the red "(This is a change ... alternate)" text to the
#     first paragraph, as three separate colored runs preceded by two spaces.

$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.End = $r.End - 1        # exclude the paragraph mark from the range
$r.Collapse(0)              # collapse to the (now real) end of the text

# Two plain spaces, uncolored, matching the existing run formatting.
$r.InsertAfter("  ")
$r.Collapse(0)

$dash = [char]0x2013

$seg1 = "(This is a change " + $dash + " Ve"
$seg2 = "rsion for branch alternate"
$seg3 = ")"

$redColor = 192   # COLORREF for RGB hex C00000 (R=0xC0,G=0x00,B=0x00 -> R + G*256 + B*65536)

$start1 = $r.Start
$r.InsertAfter($seg1)
$d.Range($start1, $start1 + $seg1.Length).Font.Color = $redColor
$r.Collapse(0)

$start2 = $r.Start
$r.InsertAfter($seg2)
$d.Range($start2, $start2 + $seg2.Length).Font.Color = $redColor
$r.Collapse(0)

$start3 = $r.Start
$r.InsertAfter($seg3)
$d.Range($start3, $start3 + $seg3.Length).Font.Color = $redColor
$r.Collapse(0)

# --- Change 2: append a new, plain/empty paragraph at the very end of the
#     document body (right before the sectPr).

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$newLastPara = $d.Paragraphs($d.Paragraphs.Count)
$newLastPara.Range.Style = "Normal"
